$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.971.67"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.365.27"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'302.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'95.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D8").Value = "'0.502"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'34.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D13").Value = "'18.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "'6.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "2.730.51"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "2.354.65"
$ws.Range("E16").Value = "  +4.80%  "
$ws.Range("D17").Value = "'0.793"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "42.921.26"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'11.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "'67.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "'235.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -5.05%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").Value = "'24.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "'32.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "'5.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "'17.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "'0.0716"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("D35").Value = "'130.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -11.83%  "
$ws.Range("D36").Value = "'1.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("E37").Value = "  +3.66%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.94%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "'21.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").Value = "1.933.47"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").Value = "'2.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "'9.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.80%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'51.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'71.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
